# "testes adicionais de probabilidade"
# Adds a column J with H^2 (acertos ao quadrado) for rows 4-8, plus the
# average and population variance of that column in rows 9-10 on
# "Planilha1". Also switches the active sheet/selection back to Planilha1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Rows 7 and 8 don't have a J cell yet in the source workbook, so give
# them the same "wrap text" style already used by J1:J6 (style index 2)
# before writing formulas into them.
$ws.Range("J7:J10").WrapText = $true

# J4 = H4^2
$ws.Range("J4").Formula = "=H4^2"

# J5:J8 share one formula (H^2), mirroring the existing shared-formula
# pattern already used in columns H and I of this sheet.
$ws.Range("J5:J8").Formula = "=H5^2"

# Summary stats below the table.
$ws.Range("J9").Formula = "=AVERAGE(J4:J8)"
$ws.Range("J10").Formula = "=VAR.P(J4:J8)"

# Make Planilha1 the active sheet/tab again, with G6 selected.
$ws.Activate() | Out-Null
$ws.Range("G6").Select() | Out-Null
